# Eleventh Commit: completed the triangle test plan file.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Developer name
$ws.Range("C3").Value = "Michael obikwere"

# Test Case 1 (__init__ - Attribute set to input values.)
$ws.Range("E7").Value = "None"
$ws.Range("F7").Value = "color: ""blue""`nside_1 : 5`nside_2 : 5`nside_3 : 7"
$ws.Range("G7").Value = "Attributes Set"

# Test Case 2 (__init__ - Exception raised when color is blank)
$ws.Range("E8").Value = "None"
$ws.Range("F8").Value = "color: """"`nside_1 : 5`nside_2 : 5`nside_3 : 7"
$ws.Range("G8").Value = "ValueError"

# Test Case 3 (__init__ - Exception raised when side_1 is not an integer.)
$ws.Range("E9").Value = "None"
$ws.Range("F9").Value = "color : ""blue""`nside_1 : ""5""`nside_2 : 5`nside_3 : 7"
$ws.Range("G9").Value = "ValueError"

# Test Case 4 (__init__ - Exception raised when side_2 is not an integer.)
$ws.Range("F10").Value = "color : ""blue""`nside_1 : 5`nside_2 : ""5""`nside_3 : 7"
$ws.Range("G10").Value = "ValueError"

# Test Case 5 (__init__ - Exception raised when side_3 is not an integer.)
$ws.Range("F11").Value = "color : ""blue""`nside_1 : 5`nside_2 : 5`nside_3 : ""7"""
$ws.Range("G11").Value = "ValueError"

# Test Case 6 (__str__ - Returns string formatted appropriately)
$ws.Range("E12").Value = "Triangle(""blue"", 5, 5, 7)"
$ws.Range("F12").Value = "None"
$ws.Range("G12").Value = "The shape color is blue.`nThis triangle has three sides with lengths of 5, 5, 7 centimeters."

# Test Case 7 (calculate_area - Returns correct calculated value.)
$ws.Range("E13").Value = "Triangle(""blue"", 5, 5, 7)"
$ws.Range("F13").Value = "None"
$ws.Range("G13").Value = "Calculated Area"

# Test Case 8 (calculate_perimeter - Returns correct calculated value.)
# Row 8's E/F/G cells still carried the plain (unwrapped/no-top-border) style
# used by the still-empty rows below; bring them in line with the other
# completed rows (same look as E7:G13) before filling in the values.
$ws.Range("E7").Copy()
$ws.Range("E14:G14").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E14").Value = "Triangle(""blue"", 5, 5, 7)"
$ws.Range("F14").Value = "None"
$ws.Range("G14").Value = "Calculated Perimeter"

# Selection cosmetic update to match final saved state
$ws.Range("F7").Select()
